# Apply the commit's row re-shuffle to the "Artfynd" sheet (rows 2-22).
# Each physical record (species observation) keeps its full set of field
# values, but the rows have been re-ordered as part of an automatic
# refresh from the source database. Columns that are identical across
# every single record (C, I, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT,
# AY) do not need to be touched - only the columns that actually carry
# per-record data are rewritten here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$records = @(
    @{Row=2; A=111896638; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575087.1320314853; R=6703393.020834555; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=3; A=111896636; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575108.85141061; R=6703418.142308297; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=4; A=111896637; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575088.0587098968; R=6703396.00058554; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=5; A=111896652; B=89183; D="LC"; E=3215; F="Rödgul trumpetsvamp"; G="Craterellus lutescens"; H="(Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575066.556649723; R=6703455.751857814; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=6; A=111896635; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575037.2974304935; R=6703389.027347369; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=7; A=111884133; B=88899; D="NT"; E=3286; F="Flattoppad klubbsvamp"; G="Clavariadelphus truncatus"; H="(Quél.) Donk"; P="Kalkberget (Kalkberget), Gstr"; Q=575059.034285416; R=6703389.477814267; AW="Patric Engfeldt"; AX="Patric Engfeldt"},
    @{Row=8; A=111884471; B=88899; D="NT"; E=3286; F="Flattoppad klubbsvamp"; G="Clavariadelphus truncatus"; H="(Quél.) Donk"; P="Kalkberget (Kalkberget), Gstr"; Q=575020.8210917887; R=6703397.074168184; AW="Patric Engfeldt"; AX="Patric Engfeldt"},
    @{Row=9; A=111883983; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kalkberget (Kalkberget), Gstr"; Q=575058.3527020445; R=6703446.206921679; AW="Patric Engfeldt"; AX="Patric Engfeldt"},
    @{Row=10; A=111896634; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575048.3395925189; R=6703452.413791304; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=11; A=111896690; B=90687; D="LC"; E=5964; F="Fjällig taggsvamp s.str."; G="Sarcodon imbricatus s.str."; H="(L.:Fr.) P.Karst."; P="Kratte masugn, Gstr"; Q=575060.2881161601; R=6703376.67477417; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=12; A=111884093; B=98535; D="LC"; E=222498; F="Blåsippa"; G="Hepatica nobilis"; H="Schreb."; P="Kopparåsen (Kopparåsen), Gstr"; Q=575065.9914513066; R=6703387.648325931; AW="Patric Engfeldt"; AX="Patric Engfeldt"},
    @{Row=13; A=111896640; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575025.3556637274; R=6703369.042946251; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=14; A=111896654; B=89183; D="LC"; E=3215; F="Rödgul trumpetsvamp"; G="Craterellus lutescens"; H="(Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575072.6962527435; R=6703421.833381963; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=15; A=111896655; B=89183; D="LC"; E=3215; F="Rödgul trumpetsvamp"; G="Craterellus lutescens"; H="(Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575104.6742508161; R=6703428.910891063; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=16; A=111896643; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575038.7114136803; R=6703416.194821274; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=17; A=111896653; B=89183; D="LC"; E=3215; F="Rödgul trumpetsvamp"; G="Craterellus lutescens"; H="(Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575075.050630242; R=6703403.625642136; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=18; A=111896639; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575089.384229039; R=6703379.745088123; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=19; A=111896644; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575036.4083237475; R=6703431.936489306; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=20; A=111896641; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575021.3626164712; R=6703370.933926445; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=21; A=111896633; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575100.4050603262; R=6703444.118284944; AW="Philipp Weiss"; AX="Philipp Weiss"},
    @{Row=22; A=111896642; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; P="Kratte masugn, Gstr"; Q=575014.1091647458; R=6703387.066676207; AW="Philipp Weiss"; AX="Philipp Weiss"}
)

foreach ($rec in $records) {
    $r = $rec.Row
    $ws.Range("A$r").Value = $rec.A
    $ws.Range("B$r").Value = $rec.B
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("E$r").Value = $rec.E
    $ws.Range("F$r").Value = $rec.F
    $ws.Range("G$r").Value = $rec.G
    $ws.Range("H$r").Value = $rec.H
    $ws.Range("P$r").Value = $rec.P
    $ws.Range("Q$r").Value = $rec.Q
    $ws.Range("R$r").Value = $rec.R
    $ws.Range("AW$r").Value = $rec.AW
    $ws.Range("AX$r").Value = $rec.AX
}

# The "Alder-Stadium" (K) and "Bestamningsmetod" (AF) helper columns are
# blank placeholder cells that travel with their record; move them to
# their new rows as well.
$kClear = @(3, 18, 20, 22)
foreach ($r in $kClear) { $ws.Range("K$r").ClearContents() }

$kAdd = @(7, 8, 9, 12)
foreach ($r in $kAdd) { $ws.Range("K$r").Value = "" }

$afClear = @(4)
foreach ($r in $afClear) { $ws.Range("AF$r").ClearContents() }

$afAdd = @(11)
foreach ($r in $afAdd) { $ws.Range("AF$r").Value = "" }
